# Auto-generated: apply numeric updates to H:N columns across multiple sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 14306533
$ws.Range("J40").Value = 35748356
$ws.Range("L40").Value = 35748356
$ws.Range("N40").Value = -35748706
$ws.Range("H80").Value = 1071.6364
$ws.Range("I80").Value = 712.5714
$ws.Range("J80").Value = 1700
$ws.Range("K80").Value = 2137.7142
$ws.Range("L80").Value = 5100
$ws.Range("M80").Value = -1139.7142
$ws.Range("N80").Value = -7096
$ws.Range("H83").Value = 1071.6364
$ws.Range("I83").Value = 712.5714
$ws.Range("J83").Value = 1700
$ws.Range("K83").Value = 6413.1426
$ws.Range("L83").Value = 15300
$ws.Range("M83").Value = -1421.1426
$ws.Range("N83").Value = -25284
$ws.Range("H98").Value = 1574.8667
$ws.Range("I98").Value = 1574.8667
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1574.8667
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -76.86670000000004
$ws.Range("N98").Value = $null
$ws.Range("H115").Value = 294.5
$ws.Range("I115").Value = 294.5
$ws.Range("K115").Value = 883.5
$ws.Range("M115").Value = 683.5
$ws.Range("H118").Value = 1302
$ws.Range("J118").Value = 300
$ws.Range("L118").Value = 900
$ws.Range("N118").Value = -4214
$ws.Range("H122").Value = 1574.8667
$ws.Range("I122").Value = 1574.8667
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4724.6001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2274.6001
$ws.Range("N122").Value = $null
$ws.Range("H129").Value = 1698.6666
$ws.Range("I129").Value = 981.2222
$ws.Range("J129").Value = 2416.111
$ws.Range("K129").Value = 2943.6666
$ws.Range("L129").Value = 7248.333
$ws.Range("M129").Value = 2056.3334
$ws.Range("N129").Value = -17248.333
$ws.Range("H132").Value = 157803.02
$ws.Range("I132").Value = 363848.03
$ws.Range("J132").Value = 22247.078
$ws.Range("K132").Value = 1091544.09
$ws.Range("L132").Value = 66741.234
$ws.Range("M132").Value = -1089014.09
$ws.Range("N132").Value = -71801.234
$ws.Range("H138").Value = 5465.8887
$ws.Range("I138").Value = 2121.8333
$ws.Range("J138").Value = 5883.896
$ws.Range("K138").Value = 6365.499899999999
$ws.Range("L138").Value = 17651.688
$ws.Range("M138").Value = -1225.499899999999
$ws.Range("N138").Value = -27931.688

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3960.9
$ws.Range("I32").Value = 2079
$ws.Range("J32").Value = 16716
$ws.Range("K32").Value = 2079
$ws.Range("L32").Value = 16716
$ws.Range("M32").Value = -1792
$ws.Range("N32").Value = -17290
$ws.Range("H45").Value = 82911
$ws.Range("I45").Value = 104341.63
$ws.Range("J45").Value = 4332
$ws.Range("K45").Value = 104341.63
$ws.Range("L45").Value = 4332
$ws.Range("M45").Value = -103964.63
$ws.Range("N45").Value = -5086
$ws.Range("H61").Value = 9620.947
$ws.Range("I61").Value = 9988.777
$ws.Range("K61").Value = 9988.777
$ws.Range("M61").Value = -9776.777
$ws.Range("H74").Value = 19232754
$ws.Range("I74").Value = 31250786
$ws.Range("K74").Value = 31250786
$ws.Range("M74").Value = -31249912
$ws.Range("H77").Value = 19232754
$ws.Range("I77").Value = 31250786
$ws.Range("K77").Value = 156253930
$ws.Range("M77").Value = -156249562
$ws.Range("H132").Value = 34205.7
$ws.Range("I132").Value = 53301.273
$ws.Range("K132").Value = 159903.819
$ws.Range("M132").Value = -157373.819
$ws.Range("H136").Value = 9620.947
$ws.Range("I136").Value = 9988.777
$ws.Range("K136").Value = 29966.331
$ws.Range("M136").Value = -27416.331
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3289.5833
$ws.Range("I20").Value = 2753.5715
$ws.Range("J20").Value = 4040
$ws.Range("K20").Value = 2753.5715
$ws.Range("L20").Value = 4040
$ws.Range("M20").Value = -2506.5715
$ws.Range("N20").Value = -4534
$ws.Range("H94").Value = 596567.9399999999
$ws.Range("I94").Value = 857043.25
$ws.Range("J94").Value = 1195.8572
$ws.Range("K94").Value = 857043.25
$ws.Range("L94").Value = 1195.8572
$ws.Range("M94").Value = -856592.25
$ws.Range("N94").Value = -2097.8572
$ws.Range("H105").Value = 1587.0938
$ws.Range("I105").Value = 1700.0588
$ws.Range("J105").Value = 1459.0667
$ws.Range("K105").Value = 1700.0588
$ws.Range("L105").Value = 1459.0667
$ws.Range("M105").Value = 46.94119999999998
$ws.Range("N105").Value = -4953.0667
$ws.Range("H118").Value = 55000
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 55000
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 55000
$ws.Range("M118").Value = $null
$ws.Range("N118").Value = -58314

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 15162706
$ws.Range("I132").Value = 17556466
$ws.Range("K132").Value = 52669398
$ws.Range("M132").Value = -52666868

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 243977.88
$ws.Range("I5").Value = 527
$ws.Range("K5").Value = 1581
$ws.Range("M5").Value = -1469
$ws.Range("H7").Value = 255
$ws.Range("I7").Value = 239.33333
$ws.Range("K7").Value = 717.99999
$ws.Range("M7").Value = -605.99999
$ws.Range("H9").Value = 898.6
$ws.Range("J9").Value = 898.5
$ws.Range("L9").Value = 2695.5
$ws.Range("N9").Value = -3143.5
$ws.Range("H135").Value = 243977.88
$ws.Range("I135").Value = 527
$ws.Range("K135").Value = 4743
$ws.Range("M135").Value = -2208
$ws.Range("H141").Value = 11146.789
$ws.Range("I141").Value = 5907.231
$ws.Range("K141").Value = 17721.693
$ws.Range("M141").Value = -12541.693

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9530298
$ws.Range("I70").Value = 15878163
$ws.Range("K70").Value = 15878163
$ws.Range("M70").Value = -15877893
$ws.Range("H73").Value = 9530298
$ws.Range("I73").Value = 15878163
$ws.Range("K73").Value = 15878163
$ws.Range("M73").Value = -15877227
$ws.Range("H97").Value = 601.90625
$ws.Range("J97").Value = 494
$ws.Range("L97").Value = 494
$ws.Range("N97").Value = -1486
$ws.Range("H102").Value = 17248944
$ws.Range("I102").Value = 23816866
$ws.Range("K102").Value = 23816866
$ws.Range("M102").Value = -23815244
$ws.Range("H122").Value = 397712.8
$ws.Range("I122").Value = 553198.1
$ws.Range("J122").Value = 8999.5
$ws.Range("K122").Value = 1659594.3
$ws.Range("L122").Value = 26998.5
$ws.Range("M122").Value = -1657144.3
$ws.Range("N122").Value = -31898.5
$ws.Range("H132").Value = 104268.85
$ws.Range("I132").Value = 137032.86
$ws.Range("K132").Value = 411098.58
$ws.Range("M132").Value = -408568.58

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 41669630
$ws.Range("I40").Value = 3284.6667
$ws.Range("K40").Value = 3284.6667
$ws.Range("M40").Value = -3148.6667
$ws.Range("H100").Value = 4762
$ws.Range("J100").Value = 4500
$ws.Range("L100").Value = 4500
$ws.Range("N100").Value = -5582
$ws.Range("H138").Value = 97714.5
$ws.Range("J138").Value = 97714.5
$ws.Range("L138").Value = 97714.5
$ws.Range("N138").Value = -107994.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 35000
$ws.Range("J49").Value = 35000
$ws.Range("L49").Value = 35000
$ws.Range("N49").Value = -35460
$ws.Range("H81").Value = 2616280.5
$ws.Range("I81").Value = 3476880.2
$ws.Range("J81").Value = 2099920.8
$ws.Range("K81").Value = 6953760.4
$ws.Range("L81").Value = 4199841.6
$ws.Range("M81").Value = -6952699.4
$ws.Range("N81").Value = -4201963.6
$ws.Range("H84").Value = 2616280.5
$ws.Range("I84").Value = 3476880.2
$ws.Range("J84").Value = 2099920.8
$ws.Range("K84").Value = 34768802
$ws.Range("L84").Value = 20999208
$ws.Range("M84").Value = -34763498
$ws.Range("N84").Value = -21009816
$ws.Range("H132").Value = 55558276
$ws.Range("I132").Value = 11112553
$ws.Range("J132").Value = 100004000
$ws.Range("K132").Value = 33337659
$ws.Range("L132").Value = 300012000
$ws.Range("M132").Value = -33335129
$ws.Range("N132").Value = -300017060
$ws.Range("H136").Value = 9164.950999999999
$ws.Range("I136").Value = 2740.6843
$ws.Range("J136").Value = 12071.167
$ws.Range("K136").Value = 8222.052899999999
$ws.Range("L136").Value = 36213.501
$ws.Range("M136").Value = -5672.052899999999
$ws.Range("N136").Value = -41313.501
$ws.Range("H137").Value = 98999.336
$ws.Range("J137").Value = 98999.336
$ws.Range("L137").Value = 98999.336
$ws.Range("N137").Value = -109199.336
